# Add latest ARS data (commit: "Add latest ARS data.")
#
# The upstream source file gained a new weekly extract (2021-03-09, ISO
# week 7) and, as a side effect of the underlying raw data for week 5
# being corrected, two existing cells feeding a couple of interpolation-
# style helper formulas changed too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected upstream values that ripple into the J169:J175 helper area ---
# J171 (=J170/J169) and J175 (=(J174-J172)*J171+J173) are formulas, so they
# recalculate automatically once their precedents are updated.
$ws.Range("J169").Value = 13
$ws.Range("J174").Value = 11.5

# --- New "Source" label row for the freshly added week ---
$ws.Range("B183").Value = "Source: 2021-03-09"

# --- New week-7 data block (age groups 0-4 .. >=80) ---
$weekData = @(
  @{ Row = 184; Age = "0-4";   NTests = 7500;   Pct = 5.6 },
  @{ Row = 185; Age = "5-14";  NTests = 11250;  Pct = 9.4 },
  @{ Row = 186; Age = "15-34"; NTests = 93750;  Pct = 5.6 },
  @{ Row = 187; Age = "35-59"; NTests = 150000; Pct = 5.5 },
  @{ Row = 188; Age = "60-79"; NTests = 84750;  Pct = 5.7 },
  @{ Row = 189; Age = ">=80";  NTests = 46500;  Pct = 7.5 }
)

foreach ($entry in $weekData) {
    $r = $entry.Row
    $ws.Range("A$r").Value = 2021
    $ws.Range("B$r").Value = 7
    $ws.Range("C$r").Value = $entry.Age
    $ws.Range("D$r").Value = $entry.NTests
    $ws.Range("E$r").Value = $entry.Pct
}

# --- Restore the saved cursor/selection position for the active pane ---
$ws.Range("H175").Select()
